# "Fruta / hortaliza, semanal" — insert a new weekly price-report row.
#
# The sheet "Hortaliza, Vega Monumental Concepción - Cebollín" gets a new
# record inserted at row 47 (shifting the former rows 47-69 down to 48-70),
# growing the used range from A1:R69 to A1:R70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 47, pushing everything from 47 downward one row.
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new weekly observation.
$ws.Range("A47").Value = 11
$ws.Range("B47").Value = "Vega Monumental Concepción"
$ws.Range("C47").Value = "Bíobío"
$ws.Range("D47").Value = 44806
$ws.Range("E47").Value = 8
$ws.Range("F47").Value = 100112037
$ws.Range("G47").Value = "Cebollín"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 220
$ws.Range("K47").Value = 4000
$ws.Range("L47").Value = 4500
$ws.Range("M47").Value = 4227
$ws.Range("N47").Value = "`$/paquete 36 unidades"
$ws.Range("O47").Value = "Región Metropolitana"
$ws.Range("P47").Value = 117
$ws.Range("Q47").Value = 36
$ws.Range("R47").Value = "Hortaliza"
